$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 previously held the text "R40"; it now holds the text "1"
# (stored as a shared string, NOT a number) while keeping the cell's
# existing style (s="23") completely untouched.
#
# A plain `.Value = "1"` assignment lets Excel auto-detect the
# digits-only string as a number, which both loses the text type and
# pulls in a new/auto style for the cell. Forcing it via NumberFormat
# "@" has a similar side effect: it permanently registers a new style
# (custom number format) in the workbook, even after being reset.
#
# Instead, stage the literal text "1" in a scratch cell using a
# formula that evaluates to text ( ="1" ), copy it, and paste-special
# only the values into B11. Pasting a text-formula's value keeps it a
# real text value (not a number) and does not touch B11's existing
# style/number format at all. The scratch cell is cleared afterwards
# so the sheet's used range/content is left exactly as before.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
